$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C) from 45190 to 45192 for all data rows (2-91)
for ($r = 2; $r -le 91; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45192
}

# Row 5 specific updates: remove "Leptoporus erubescens" finding, adjust related counts
$ws.Cells.Item(5, 10).Value2 = 3   # J5: NT count 4 -> 3
$ws.Cells.Item(5, 15).Value2 = 3   # O5: Rödlistade count 4 -> 3
$ws.Cells.Item(5, 17).Value2 = 16  # Q5: Alla arter count 17 -> 16

$r5Text = "Barrviolspindling`r`nSpillkråka`r`nVedtrappmossa`r`nBronshjon`r`nDropptaggsvamp`r`nFällmossa`r`nGrön sköldmossa`r`nGuldlockmossa`r`nKornknutmossa`r`nRödgul trumpetsvamp`r`nStubbspretmossa`r`nSvavelriska`r`nSårläka`r`nVågbandad barkbock`r`nBlåsippa`r`nRevlummer"
$ws.Cells.Item(5, 18).Value = $r5Text
